$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6193925142288208
$ws.Range("B1").Value = 1.410785675048828
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.382034301757812
$ws.Range("E1").Value = 1.382242560386658
